# Re-style the three tables (slides 14, 15, 16) that use the old
# "Table_0" table style with the built-in "No Style, No Grid" table
# style, matching the GUID swap captured in the commit.
$p = $ppt.ActivePresentation

$oldStyleId = "{64D831A6-858E-4A50-978F-8030A17132FB}"
$newStyleId = "{E4864D30-3C3A-48F8-955D-05C91EBF22F3}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
